# The template's logos were re-exported under new file names:
#   - the BTEC logo picture in the "first page" header
#     (currently stored as image2.jpg) should become image1.jpg
#   - the Pearson/Edexcel logo pictures in the footers
#     (currently stored as image1.png) should become image2.png
#
# Walk the (single) section's headers and footers, find each inline
# picture by its description, and rename it.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}

Write-Output "Renamed header logo(s)."

for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}

Write-Output "Renamed footer logo(s)."
